# Apply "完成日内简单调度" edits to PrivateCar_example.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data rows 2-11: column D (P): 3 -> 2, column E (P_variance): 0.5 -> 0.3,
# column Q (Car_area_end): 3 -> 4
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 4).Value = 2     # D column
    $ws.Cells.Item($r, 5).Value = 0.3   # E column
    $ws.Cells.Item($r, 17).Value = 4    # Q column
}

# Update the active cell selection on the sheet
$ws.Activate()
$ws.Range("P19").Select()
